# Weekly update: insert two new daily price records (Betarraga, Vega Central
# Mapocho de Santiago) ahead of the existing history, shifting the prior
# rows (548 onward) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 548; this shifts the old
# rows 548:593 down to 550:595 and carries their formatting (e.g. the
# date-style on column D) the same way Excel's own Rows.Insert does.
$ws.Range("A548:A549").EntireRow.Insert()

# New row 548 - "Primera" quality record for 2022-08-13 (serial 44783)
$ws.Range("A548").Value = 9
$ws.Range("B548").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C548").Value = "Metropolitana"
$ws.Range("D548").Value = 44783
$ws.Range("E548").Value = 13
$ws.Range("F548").Value = 100114014
$ws.Range("G548").Value = "Betarraga"
$ws.Range("H548").Value = "Sin especificar"
$ws.Range("I548").Value = "Primera"
$ws.Range("J548").Value = 4300
$ws.Range("K548").Value = 180
$ws.Range("L548").Value = 190
$ws.Range("M548").Value = 185
$ws.Range("N548").Value = "$/unidad"
$ws.Range("O548").Value = "Región Metropolitana"
$ws.Range("P548").Value = 185
$ws.Range("Q548").Value = 1
$ws.Range("R548").Value = "Hortaliza"

# New row 549 - "Segunda" quality record for 2022-08-13 (serial 44783)
$ws.Range("A549").Value = 9
$ws.Range("B549").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C549").Value = "Metropolitana"
$ws.Range("D549").Value = 44783
$ws.Range("E549").Value = 13
$ws.Range("F549").Value = 100114014
$ws.Range("G549").Value = "Betarraga"
$ws.Range("H549").Value = "Sin especificar"
$ws.Range("I549").Value = "Segunda"
$ws.Range("J549").Value = 7900
$ws.Range("K549").Value = 150
$ws.Range("L549").Value = 160
$ws.Range("M549").Value = 155
$ws.Range("N549").Value = "$/unidad"
$ws.Range("O549").Value = "Región Metropolitana"
$ws.Range("P549").Value = 155
$ws.Range("Q549").Value = 1
$ws.Range("R549").Value = "Hortaliza"
